# Fruta / hortaliza, semanal
# Update the weekly Caqui price records (rows 2,3,5,6,7) with refreshed
# values for Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio,
# Origen and Precio $/Kg. Row 4 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44355
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1139

# Row 3
$ws.Range("D3").Value = 44301
$ws.Range("K3").Value = "Hachiya"
$ws.Range("M3").Value = 250
$ws.Range("R3").Value = "Región de O'Higgins"

# Row 5
$ws.Range("D5").Value = 44342
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("S5").Value = 1361

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1194

# Row 7
$ws.Range("D7").Value = 44699
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("S7").Value = 1639
